$d = $word.ActiveDocument
$wdReplaceAll = 2

# --- Locate the two Non-Functional-Requirements paragraphs we need to rework ---
$r1Index = $null
$r2Index = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "R1: The program must be written*") { $r1Index = $i }
    if ($t -like "R2: The program must be finished*") { $r2Index = $i }
}

# --- Rewrite the old "R1" paragraph (C#/Golang/Socket/SQLite) as the new R1 ---
# (old) R1: The program must be written using C# for the front-end, Golang/Python
#       for the back-end, Socket for network communication and SQLite for the database.
# (new) R1: The program must implement the client/server model.
$p1 = $d.Paragraphs.Item($r1Index)
$oldR1Text = "R1: The program must be written using C# for the front-end, Golang/Python for the back-end, Socket for network communication and SQLite for the database."
$p1.Range.Find.Execute($oldR1Text, $true, $false, $false, $false, $false, $true, 1, $false,
    "R1: The program must implement the client/server model.", $wdReplaceAll)

# --- Insert two brand-new list items between the (new) R1 and the old R2 paragraph ---
# These become R2 and R3; the old R2 paragraph itself will be turned into R4 below.
$p1Again = $d.Paragraphs.Item($r1Index)
$p1Again.Range.InsertParagraphAfter()
$newR2 = $d.Paragraphs.Item($r1Index + 1)
$newR2.Range.Text = "R2: The program must use a database."

$newR2Again = $d.Paragraphs.Item($r1Index + 1)
$newR2Again.Range.InsertParagraphAfter()
$newR3 = $d.Paragraphs.Item($r1Index + 2)
$newR3.Range.Text = "R3: The program must be written using C# for the front-end, Golang/Python for the back-end, Socket for network communication and SQLite for the database."

# --- The paragraph that used to be "R2: ... fall semester." becomes "R4: ..." ---
$r2Index = $r1Index + 3
$p4 = $d.Paragraphs.Item($r2Index)
$oldR2Text = "R2: The program must be finished within the current fall semester."
$p4.Range.Find.Execute($oldR2Text, $true, $false, $false, $false, $false, $true, 1, $false,
    "R4: The program must be finished within the current fall semester.", $wdReplaceAll)

# --- Insert a brand-new R5 list item right after it ---
$p4Again = $d.Paragraphs.Item($r2Index)
$p4Again.Range.InsertParagraphAfter()
$newR5 = $d.Paragraphs.Item($r2Index + 1)
$newR5.Range.Text = "R5: The program is limited to managing reading material."
